$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 8.783898333333333
$ws.Cells.Item(2, 8).Value = 26.351695
$ws.Cells.Item(2, 9).Value = 0.09847125088802929
$ws.Cells.Item(2, 10).Value = 0.09847125088802929
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 3.795192333333334
$ws.Cells.Item(2, 14).Value = 11.385577
$ws.Cells.Item(2, 15).Value = 0.01044213755712683
$ws.Cells.Item(2, 16).Value = 0.01044213755712683
$ws.Cells.Item(2, 17).Value = 33.33658361144612
$ws.Cells.Item(2, 18).Value = 300.029252503015
$ws.Cells.Item(2, 19).Value = 0.00102825034719515
$ws.Cells.Item(2, 20).Value = 0.00102825034719515

$ws.Cells.Item(3, 7).Value = 8.783898333333333
$ws.Cells.Item(3, 8).Value = 26.351695
$ws.Cells.Item(3, 9).Value = 0.09847125088802929
$ws.Cells.Item(3, 10).Value = 0.09847125088802929
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 243.3763986666667
$ws.Cells.Item(3, 14).Value = 730.1291960000001
$ws.Cells.Item(3, 15).Value = 0.6696287328350964
$ws.Cells.Item(3, 16).Value = 0.6696287328350964
$ws.Cells.Item(3, 17).Value = 2137.793542620802
$ws.Cells.Item(3, 18).Value = 19240.14188358722
$ws.Cells.Item(3, 19).Value = 0.06593917895283791
$ws.Cells.Item(3, 20).Value = 0.06593917895283791

$ws.Cells.Item(4, 7).Value = 8.783898333333333
$ws.Cells.Item(4, 8).Value = 26.351695
$ws.Cells.Item(4, 9).Value = 0.09847125088802929
$ws.Cells.Item(4, 10).Value = 0.09847125088802929
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 29.801371
$ws.Cells.Item(4, 14).Value = 89.404113
$ws.Cells.Item(4, 15).Value = 0.08199584844219236
$ws.Cells.Item(4, 16).Value = 0.08199584844219235
$ws.Cells.Item(4, 17).Value = 261.7722130579483
$ws.Cells.Item(4, 18).Value = 2355.949917521535
$ws.Cells.Item(4, 19).Value = 0.008074233763727949
$ws.Cells.Item(4, 20).Value = 0.008074233763727947

$ws.Cells.Item(5, 7).Value = 8.783898333333333
$ws.Cells.Item(5, 8).Value = 26.351695
$ws.Cells.Item(5, 9).Value = 0.09847125088802929
$ws.Cells.Item(5, 10).Value = 0.09847125088802929
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 86.47679266666667
$ws.Cells.Item(5, 14).Value = 259.430378
$ws.Cells.Item(5, 15).Value = 0.2379332811655844
$ws.Cells.Item(5, 16).Value = 0.2379332811655844
$ws.Cells.Item(5, 17).Value = 759.6033549767455
$ws.Cells.Item(5, 18).Value = 6836.430194790711
$ws.Cells.Item(5, 19).Value = 0.02342958782426827
$ws.Cells.Item(5, 20).Value = 0.02342958782426827

$ws.Cells.Item(6, 7).Value = 45.41653666666667
$ws.Cells.Item(6, 8).Value = 136.24961
$ws.Cells.Item(6, 9).Value = 0.5091387681022471
$ws.Cells.Item(6, 10).Value = 0.5091387681022471
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 3.795192333333334
$ws.Cells.Item(6, 14).Value = 11.385577
$ws.Cells.Item(6, 15).Value = 0.01044213755712683
$ws.Cells.Item(6, 16).Value = 0.01044213755712683
$ws.Cells.Item(6, 17).Value = 172.3644917638856
$ws.Cells.Item(6, 18).Value = 1551.28042587497
$ws.Cells.Item(6, 19).Value = 0.005316497052189764
$ws.Cells.Item(6, 20).Value = 0.005316497052189764

$ws.Cells.Item(7, 7).Value = 45.41653666666667
$ws.Cells.Item(7, 8).Value = 136.24961
$ws.Cells.Item(7, 9).Value = 0.5091387681022471
$ws.Cells.Item(7, 10).Value = 0.5091387681022471
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 243.3763986666667
$ws.Cells.Item(7, 14).Value = 730.1291960000001
$ws.Cells.Item(7, 15).Value = 0.6696287328350964
$ws.Cells.Item(7, 16).Value = 0.6696287328350964
$ws.Cells.Item(7, 17).Value = 11053.31313384595
$ws.Cells.Item(7, 18).Value = 99479.81820461359
$ws.Cells.Item(7, 19).Value = 0.3409339481215297
$ws.Cells.Item(7, 20).Value = 0.3409339481215297

$ws.Cells.Item(8, 7).Value = 45.41653666666667
$ws.Cells.Item(8, 8).Value = 136.24961
$ws.Cells.Item(8, 9).Value = 0.5091387681022471
$ws.Cells.Item(8, 10).Value = 0.5091387681022471
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 29.801371
$ws.Cells.Item(8, 14).Value = 89.404113
$ws.Cells.Item(8, 15).Value = 0.08199584844219236
$ws.Cells.Item(8, 16).Value = 0.08199584844219235
$ws.Cells.Item(8, 17).Value = 1353.475058738437
$ws.Cells.Item(8, 18).Value = 12181.27552864593
$ws.Cells.Item(8, 19).Value = 0.04174726526535637
$ws.Cells.Item(8, 20).Value = 0.04174726526535637

$ws.Cells.Item(9, 7).Value = 45.41653666666667
$ws.Cells.Item(9, 8).Value = 136.24961
$ws.Cells.Item(9, 9).Value = 0.5091387681022471
$ws.Cells.Item(9, 10).Value = 0.5091387681022471
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 86.47679266666667
$ws.Cells.Item(9, 14).Value = 259.430378
$ws.Cells.Item(9, 15).Value = 0.2379332811655844
$ws.Cells.Item(9, 16).Value = 0.2379332811655844
$ws.Cells.Item(9, 17).Value = 3927.476424961398
$ws.Cells.Item(9, 18).Value = 35347.28782465259
$ws.Cells.Item(9, 19).Value = 0.1211410576631712
$ws.Cells.Item(9, 20).Value = 0.1211410576631712

$ws.Cells.Item(10, 7).Value = 25.203909
$ws.Cells.Item(10, 8).Value = 75.611727
$ws.Cells.Item(10, 9).Value = 0.2825465815194877
$ws.Cells.Item(10, 10).Value = 0.2825465815194877
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 3.795192333333334
$ws.Cells.Item(10, 14).Value = 11.385577
$ws.Cells.Item(10, 15).Value = 0.01044213755712683
$ws.Cells.Item(10, 16).Value = 0.01044213755712683
$ws.Cells.Item(10, 17).Value = 95.65368220683101
$ws.Cells.Item(10, 18).Value = 860.8831398614791
$ws.Cells.Item(10, 19).Value = 0.002950390270522441
$ws.Cells.Item(10, 20).Value = 0.002950390270522441

$ws.Cells.Item(11, 7).Value = 25.203909
$ws.Cells.Item(11, 8).Value = 75.611727
$ws.Cells.Item(11, 9).Value = 0.2825465815194877
$ws.Cells.Item(11, 10).Value = 0.2825465815194877
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 243.3763986666667
$ws.Cells.Item(11, 14).Value = 730.1291960000001
$ws.Cells.Item(11, 15).Value = 0.6696287328350964
$ws.Cells.Item(11, 16).Value = 0.6696287328350964
$ws.Cells.Item(11, 17).Value = 6134.036604742389
$ws.Cells.Item(11, 18).Value = 55206.3294426815
$ws.Cells.Item(11, 19).Value = 0.1892013093497829
$ws.Cells.Item(11, 20).Value = 0.1892013093497829

$ws.Cells.Item(12, 7).Value = 25.203909
$ws.Cells.Item(12, 8).Value = 75.611727
$ws.Cells.Item(12, 9).Value = 0.2825465815194877
$ws.Cells.Item(12, 10).Value = 0.2825465815194877
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 29.801371
$ws.Cells.Item(12, 14).Value = 89.404113
$ws.Cells.Item(12, 15).Value = 0.08199584844219236
$ws.Cells.Item(12, 16).Value = 0.08199584844219235
$ws.Cells.Item(12, 17).Value = 751.111042759239
$ws.Cells.Item(12, 18).Value = 6759.999384833151
$ws.Cells.Item(12, 19).Value = 0.02316764667613147
$ws.Cells.Item(12, 20).Value = 0.02316764667613146

$ws.Cells.Item(13, 7).Value = 25.203909
$ws.Cells.Item(13, 8).Value = 75.611727
$ws.Cells.Item(13, 9).Value = 0.2825465815194877
$ws.Cells.Item(13, 10).Value = 0.2825465815194877
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 86.47679266666667
$ws.Cells.Item(13, 14).Value = 259.430378
$ws.Cells.Item(13, 15).Value = 0.2379332811655844
$ws.Cells.Item(13, 16).Value = 0.2379332811655844
$ws.Cells.Item(13, 17).Value = 2179.553212982534
$ws.Cells.Item(13, 18).Value = 19615.97891684281
$ws.Cells.Item(13, 19).Value = 0.06722723522305099
$ws.Cells.Item(13, 20).Value = 0.06722723522305099

$ws.Cells.Item(14, 7).Value = 9.798324333333333
$ws.Cells.Item(14, 8).Value = 29.394973
$ws.Cells.Item(14, 9).Value = 0.1098433994902357
$ws.Cells.Item(14, 10).Value = 0.1098433994902357
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 3.795192333333334
$ws.Cells.Item(14, 14).Value = 11.385577
$ws.Cells.Item(14, 15).Value = 0.01044213755712683
$ws.Cells.Item(14, 16).Value = 0.01044213755712683
$ws.Cells.Item(14, 17).Value = 37.18652538938012
$ws.Cells.Item(14, 18).Value = 334.678728504421
$ws.Cells.Item(14, 19).Value = 0.001146999887219477
$ws.Cells.Item(14, 20).Value = 0.001146999887219477

$ws.Cells.Item(15, 7).Value = 9.798324333333333
$ws.Cells.Item(15, 8).Value = 29.394973
$ws.Cells.Item(15, 9).Value = 0.1098433994902357
$ws.Cells.Item(15, 10).Value = 0.1098433994902357
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 243.3763986666667
$ws.Cells.Item(15, 14).Value = 730.1291960000001
$ws.Cells.Item(15, 15).Value = 0.6696287328350964
$ws.Cells.Item(15, 16).Value = 0.6696287328350964
$ws.Cells.Item(15, 17).Value = 2384.680889214635
$ws.Cells.Item(15, 18).Value = 21462.12800293171
$ws.Cells.Item(15, 19).Value = 0.07355429641094581
$ws.Cells.Item(15, 20).Value = 0.07355429641094581

$ws.Cells.Item(16, 7).Value = 9.798324333333333
$ws.Cells.Item(16, 8).Value = 29.394973
$ws.Cells.Item(16, 9).Value = 0.1098433994902357
$ws.Cells.Item(16, 10).Value = 0.1098433994902357
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 29.801371
$ws.Cells.Item(16, 14).Value = 89.404113
$ws.Cells.Item(16, 15).Value = 0.08199584844219236
$ws.Cells.Item(16, 16).Value = 0.08199584844219235
$ws.Cells.Item(16, 17).Value = 292.0034986359943
$ws.Cells.Item(16, 18).Value = 2628.031487723949
$ws.Cells.Item(16, 19).Value = 0.009006702736976556
$ws.Cells.Item(16, 20).Value = 0.009006702736976555

$ws.Cells.Item(17, 7).Value = 9.798324333333333
$ws.Cells.Item(17, 8).Value = 29.394973
$ws.Cells.Item(17, 9).Value = 0.1098433994902357
$ws.Cells.Item(17, 10).Value = 0.1098433994902357
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 86.47679266666667
$ws.Cells.Item(17, 14).Value = 259.430378
$ws.Cells.Item(17, 15).Value = 0.2379332811655844
$ws.Cells.Item(17, 16).Value = 0.2379332811655844
$ws.Cells.Item(17, 17).Value = 847.3276618544215
$ws.Cells.Item(17, 18).Value = 7625.948956689795
$ws.Cells.Item(17, 19).Value = 0.02613540045509386
$ws.Cells.Item(17, 20).Value = 0.02613540045509386

